# Updated Global Glider Cal and Ingest sheets:
#  - Changed Cal scattering angle (CC_scattering_angle) value to 140
#  - Changed Cal angular resolution (CC_angular_resolution) value to 1.13

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

# Row 2: CC_scattering_angle -> 140
$ws.Range("F2").Value = 140

# Row 4: CC_angular_resolution -> 1.13
$ws.Range("F4").Value = 1.13
